$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.100.85'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '2.363.50'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '502.76'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.56'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.545'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').Value = '2.365.81'
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0985'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.80'
$ws.Range('E12').Value = '  +3.96%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.783.14'
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').Value = '56.040.25'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.39'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '2.382.85'
$ws.Range('E18').Value = '  -2.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.01'
$ws.Range('E19').Value = '  -3.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.02'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '307.81'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.00'
$ws.Range('E24').Value = '  +0.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.147'
$ws.Range('E27').Value = '  -5.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.24'
$ws.Range('E28').Value = '  -4.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.58'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').Value = '0.0₃0710'
$ws.Range('E30').Value = '  -3.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.65'
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.79'
$ws.Range('E33').Value = '  -5.51%  '
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.64'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.16'
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('E38').Value = '  -4.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.28'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.799'
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('E41').Value = '  -5.82%  '
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '129.20'
$ws.Range('E43').Value = '  -4.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.68'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.562'
$ws.Range('E45').Value = '  -2.06%  '
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '238.13'
$ws.Range('E47').Value = '  -6.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0482'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0207'
$ws.Range('E49').Value = '  -3.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.01'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('E51').Value = '  -1.24%  '
